# The "3.1.3 日志配置" heading paragraph accidentally ends with a duplicated,
# stray "配置" run placed right after the "_GoBack" bookmark
# (visible heading text is "3.1.3 日志配置配置" instead of "3.1.3 日志配置").
# Remove that trailing duplicate run while leaving the bookmark and the rest
# of the paragraph untouched.

$d = $word.ActiveDocument

$bookmarkName = "_GoBack"
$deleted = $false

if ($d.Bookmarks.Exists($bookmarkName)) {
    $bm = $d.Bookmarks.Item($bookmarkName)
    $bmEnd = $bm.End

    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        $pStart = $p.Range.Start
        $pEnd = $p.Range.End

        if ($bmEnd -ge $pStart -and $bmEnd -le $pEnd) {
            # Range spanning from just after the bookmark to just before the
            # paragraph mark (the trailing run we want to drop).
            $tailRange = $d.Range($bmEnd, $pEnd - 1)

            if ($tailRange.Text -eq "配置") {
                $tailRange.Delete()
                $deleted = $true
            }
            break
        }
    }
}

if (-not $deleted) {
    # Fallback: locate the unique duplicated heading text directly.
    $found = $d.Content.Find.Execute("日志配置配置", $true, $false, $false,
                                      $false, $false, $true, 1, $false,
                                      "日志配置", 2)
}
